$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# HOUR_APPR_PROCESS_START column (V) currently holds plain numbers (hours).
# Convert each value to a string and concatenate ':00:00' to build a time-like
# text representation, e.g. 13 -> "13:00:00".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 22).End(-4162).Row  # xlUp = -4162, col 22 = V

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 22)
    $hour = $cell.Value2
    if ($null -ne $hour) {
        $hourStr = [string]([int]$hour)
        $cell.Value2 = $hourStr + ":00:00"
    }
}
